$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: replace French "Foreigner" row with Hindi "Foreigner" row (hin/FR/परदेशी/TRUE)
$ws.Range("A3").Value = "hin"
$ws.Range("B3").Value = "FR"
$ws.Range("C3").Value = "परदेशी"

# Row 4: replace old French "Non-Foreigner" row with English "Non-Foreigner" row (eng/NFR/Non-Foreigner/TRUE)
$ws.Range("A4").Value = "eng"
$ws.Range("B4").Value = "NFR"
$ws.Range("C4").Value = "Non-Foreigner"

# Row 5: replace old French "Non-étranger" row with Hindi "Non-Foreigner" row (hin/NFR/गैर विदेशी/TRUE)
$ws.Range("A5").Value = "hin"
$ws.Range("B5").Value = "NFR"
$ws.Range("C5").Value = "गैर विदेशी"

# Remove the extra per-cell font style (s="4") from the data rows; revert to default style
$ws.Range("A2:C5").Style = "Normal"

# Drop the Arabic rows (6 & 7) entirely, leaving the rows present but blank
$ws.Rows("6:7").Delete()
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").Style = "Normal"

# Match the final selected cell recorded in the sheet
$null = $ws.Range("G11").Select()

Write-Host "done"
